# [BI-1613] Update TAF to include term type
# Adds a new "Term Type" header column (R) to the Template sheet, matching
# the header formatting already used for the other scale-related headers
# (bold font, thin border, wrapped text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# New header cell in column R, row 1
$rCell = $ws.Range("R1")
$rCell.Value = "Term Type"

# Header formatting: bold 11pt, wrap text, thin border around the cell
$rCell.Font.Bold = $true
$rCell.Font.Size = 11
$rCell.WrapText = $true
$rCell.Borders.LineStyle = 1
$rCell.Borders.Weight = 2

# Row 1 grows to fit the now-wrapped header text
$ws.Rows.Item(1).RowHeight = 32

# Keep the active selection on the newly added header cell, matching the
# post-edit workbook view.
$null = $rCell.Select()
